$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D2:D51) is treated as text so numeric-looking
# values (e.g. "1.00", "6.22") are not auto-converted to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '67.400.69'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '3.772.67'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '591.70'
$ws.Range('E5').Value = '  -3.54%  '
$ws.Range('D6').Value = '171.71'
$ws.Range('E6').Value = '  -3.85%  '
$ws.Range('D7').Value = '3.767.52'
$ws.Range('E7').Value = '  +1.09%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('E10').Value = '  -4.32%  '
$ws.Range('D11').Value = '6.22'
$ws.Range('E11').Value = '  -4.86%  '
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  -4.49%  '
$ws.Range('D13').Value = '37.78'
$ws.Range('E13').Value = '  -5.04%  '
$ws.Range('E14').Value = '  -4.01%  '
$ws.Range('D15').Value = '4.399.31'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '3.768.07'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').Value = '67.464.09'
$ws.Range('E17').Value = '  -3.17%  '
$ws.Range('E18').Value = '  -4.73%  '
$ws.Range('D19').Value = '7.08'
$ws.Range('E19').Value = '  -5.18%  '
$ws.Range('D20').Value = '16.01'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').Value = '486.34'
$ws.Range('E21').Value = '  -2.98%  '
$ws.Range('D22').Value = '9.13'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '0.719'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '83.96'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  -10.17%  '
$ws.Range('D26').Value = '0.0000139'
$ws.Range('E26').Value = '  +2.39%  '
$ws.Range('D27').Value = '12.15'
$ws.Range('E27').Value = '  -5.89%  '
$ws.Range('E28').Value = '  -11.24%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '2.90'
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').Value = '2.39'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').Value = '32.29'
$ws.Range('E32').Value = '  +6.61%  '
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  -4.01%  '
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('D38').Value = '5.71'
$ws.Range('E38').Value = '  -6.35%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.322'
$ws.Range('E39').Value = '  -8.23%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '448.98'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('D41').Value = '48.78'
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('E43').Value = '  -7.05%  '
$ws.Range('D44').Value = '8.23'
$ws.Range('E44').Value = '  -3.53%  '
$ws.Range('D45').Value = '41.34'
$ws.Range('E45').Value = '  -9.89%  '
$ws.Range('D46').Value = '2.824.49'
$ws.Range('E46').Value = '  -4.35%  '
$ws.Range('D47').Value = '139.84'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('D50').Value = '25.78'
$ws.Range('E50').Value = '  -4.91%  '
$ws.Range('D51').Value = '23.08'
$ws.Range('E51').Value = '  +7.18%  '

# Remove the temporary text formatting so the cell style matches the
# original (no explicit style index on these cells).
$priceRange.ClearFormats()
